# Apply portfolio_obligacje.xlsx updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the date column keeps storing plain text values (not auto-converted
# to Excel date serials) by pre-formatting the cells as Text before assignment.
$ws.Range("E2:E9").NumberFormat = "@"

# Row 2: ARH0227 -> BST0726, add typ dzialalnosci, wolumen 1 -> 4
$ws.Range("A2").Value = "BST0726"
$ws.Range("B2").Value = "Windykacja i zarządzanie wierzytelnościami"
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = "2025-01-05"

# Row 3: ANW1126 -> CAV0927, add typ dzialalnosci (C/D/E unchanged)
$ws.Range("A3").Value = "CAV0927"
$ws.Range("B3").Value = "Deweloper nieruchomości komercyjnych"
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "2025-01-05"

# Row 4: AOW1225 -> CAV0927, add typ dzialalnosci, cena 100 -> 99.8, wolumen 1 -> 2
$ws.Range("A4").Value = "CAV0927"
$ws.Range("B4").Value = "Deweloper nieruchomości komercyjnych"
$ws.Range("C4").Value = 99.8
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = "2025-01-05"

# Row 5: only data_dodania_obligacji changes
$ws.Range("E5").Value = "2025-01-06"

# Row 6: BST0726 -> ACH0427, typ dzialalnosci changes, cena 100 -> 103, wolumen 4 -> 1
$ws.Range("A6").Value = "ACH0427"
$ws.Range("B6").Value = "Deweloper nieruchomości"
$ws.Range("C6").Value = 103
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = "2025-01-06"

# Row 7: CAV0927 -> ACH1125, typ dzialalnosci changes, cena 100 -> 101, wolumen 1 -> 5
$ws.Range("A7").Value = "ACH1125"
$ws.Range("B7").Value = "Deweloper nieruchomości"
$ws.Range("C7").Value = 101
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = "2025-01-06"

# Row 8: CAV0927 -> BST0226, typ dzialalnosci changes, cena 99.8 -> 100, wolumen 1 -> 6
$ws.Range("A8").Value = "BST0226"
$ws.Range("B8").Value = "Windykacja i zarządzanie wierzytelnościami"
$ws.Range("C8").Value = 100
$ws.Range("D8").Value = 6
$ws.Range("E8").Value = "2025-01-06"

# Row 9: CAV0927 -> BST0226, typ dzialalnosci changes, cena 99.8 -> 100, wolumen 2 -> 6
$ws.Range("A9").Value = "BST0226"
$ws.Range("B9").Value = "Windykacja i zarządzanie wierzytelnościami"
$ws.Range("C9").Value = 100
$ws.Range("D9").Value = 6
$ws.Range("E9").Value = "2025-01-06"

$wb.Save()
